$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = 0
$r1 = $ws.Range("B1")
$r1.Value = 0

# A2 = 0
$r2 = $ws.Range("A2")
$r2.Value = 0

# B2 = "disconnected_elements" (plain text, default style)
$ws.Range("B2").Value = "disconnected_elements"

# Apply the bold + centered + thin-bordered style to B1 first
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2

# Copy B1's format onto A2 so both share the same style entry
$r1.Copy()
$r2.PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
